# Atualização de bases das ligas, do dia: 14-04-2024 às 15:12
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Swap the data (columns B..AC) between pairs of rows. Column A
#    (the running "id" index) stays where it is; every other column
#    value is exchanged between the two rows of each pair.
# ---------------------------------------------------------------
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Swap-Rows($rowA, $rowB) {
    foreach ($col in $cols) {
        $cellA = $ws.Range("$col$rowA")
        $cellB = $ws.Range("$col$rowB")
        $vA = $cellA.Value2
        $vB = $cellB.Value2
        $cellA.Value2 = $vB
        $cellB.Value2 = $vA
    }
}

Swap-Rows 29 30
Swap-Rows 36 37
Swap-Rows 87 88
Swap-Rows 111 112

# ---------------------------------------------------------------
# 2) Append the two new match rows (152 and 153) at the bottom of
#    the sheet, copying the cell formatting (styles) of columns A
#    and E from the last existing row (151) so number formats /
#    borders / bold match the rest of the table.
# ---------------------------------------------------------------
$ws.Range("A151").Copy($ws.Range("A152"))
$ws.Range("E151").Copy($ws.Range("E152"))
$ws.Range("A151").Copy($ws.Range("A153"))
$ws.Range("E151").Copy($ws.Range("E153"))

# Row 152
$ws.Range("A152").Value2 = 150
$ws.Range("B152").Value2 = 7952749
$ws.Range("C152").Value2 = "Bosnia Herzegovina Premier Liga"
$ws.Range("D152").Value2 = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E152").Value2 = 45394.4375
$ws.Range("F152").Value2 = "FK Tuzla City"
$ws.Range("G152").Value2 = "Siroki Brijeg"
$ws.Range("H152").Value2 = 3
$ws.Range("I152").Value2 = 1
$ws.Range("J152").Value2 = "H"
$ws.Range("K152").Value2 = 1.727
$ws.Range("L152").Value2 = 3.4
$ws.Range("M152").Value2 = 4.333
$ws.Range("N152").Value2 = 1.7
$ws.Range("O152").Value2 = 3.25
$ws.Range("P152").Value2 = 4.5
$ws.Range("Q152").Value2 = -0.75
$ws.Range("R152").Value2 = 1.975
$ws.Range("S152").Value2 = 1.825
$ws.Range("T152").Value2 = 2.25
$ws.Range("U152").Value2 = 1.975
$ws.Range("V152").Value2 = 1.825
$ws.Range("W152").Value2 = 0.7
$ws.Range("X152").Value2 = -1
$ws.Range("Y152").Value2 = -1
$ws.Range("Z152").Value2 = 0.9750000000000001
$ws.Range("AA152").Value2 = -1
$ws.Range("AB152").Value2 = 0.9750000000000001
$ws.Range("AC152").Value2 = -1

# Row 153
$ws.Range("A153").Value2 = 151
$ws.Range("B153").Value2 = 7952746
$ws.Range("C153").Value2 = "Bosnia Herzegovina Premier Liga"
$ws.Range("D153").Value2 = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E153").Value2 = 45394.64583333334
$ws.Range("F153").Value2 = "Borac Banja Luka"
$ws.Range("G153").Value2 = "Sloga"
$ws.Range("H153").Value2 = 1
$ws.Range("I153").Value2 = 0
$ws.Range("J153").Value2 = "H"
$ws.Range("K153").Value2 = 1.3
$ws.Range("L153").Value2 = 5.5
$ws.Range("M153").Value2 = 6
$ws.Range("N153").Value2 = 1.09
$ws.Range("O153").Value2 = 9.5
$ws.Range("P153").Value2 = 19
$ws.Range("Q153").Value2 = -2.5
$ws.Range("R153").Value2 = 1.95
$ws.Range("S153").Value2 = 1.85
$ws.Range("T153").Value2 = 3.25
$ws.Range("U153").Value2 = 1.875
$ws.Range("V153").Value2 = 1.925
$ws.Range("W153").Value2 = 0.09000000000000008
$ws.Range("X153").Value2 = -1
$ws.Range("Y153").Value2 = -1
$ws.Range("Z153").Value2 = -1
$ws.Range("AA153").Value2 = 0.8500000000000001
$ws.Range("AB153").Value2 = -1
$ws.Range("AC153").Value2 = 0.925
